$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 74: new row with the data that used to be in row 73 (original values)
$ws.Range("A74").Value = 1
$ws.Range("B74").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C74").Value = "Arica y Parinacota"
$ws.Range("D74").Value = 44592
$ws.Range("E74").Value = 15
$ws.Range("F74").Value = 100112040
$ws.Range("G74").Value = "Cilantro"
$ws.Range("H74").Value = "Sin especificar"
$ws.Range("I74").Value = "Primera"
$ws.Range("J74").Value = 300
$ws.Range("K74").Value = 1500
$ws.Range("L74").Value = 2000
$ws.Range("M74").Value = 1750
$ws.Range("N74").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O74").Value = "Región de Arica y Parinacota"
$ws.Range("P74").Value = 875
$ws.Range("Q74").Value = 2
$ws.Range("R74").Value = "Hortaliza"

# Copy the date cell style (s="2") from D73 to D74
$ws.Range("D73").Copy()
$ws.Range("D74").PasteSpecial(-4122) | Out-Null

# Row 73: updated with new values (new date + updated prices)
$ws.Range("D73").Value = 44656
$ws.Range("K73").Value = 1000
$ws.Range("L73").Value = 1500
$ws.Range("M73").Value = 1250
$ws.Range("P73").Value = 625
